# "Generate Report for Handback"
#
# The localization hand-off for both e2e test files has come back from the
# translators, so the status report is regenerated to reflect that the
# zh-cn and de-de targets are now handed back / in sync with en-US:
#   - Every "Status" cell (the Overview sheet's per-language columns, and
#     the Status column on each language's detail sheet) flips from
#     "Ready for handoff" to "Handed back: in sync with en-US".
#   - zh-cn / de-de detail sheets: the "Latest Target File" and
#     "Latest Handback File" columns are now populated (target file is a
#     hyperlink, same as the source-file link in column A), and the
#     "Latest Handback DateTime" is stamped with the handback timestamp.
#   - A few columns are widened so the new/longer values aren't clipped.

$wb = $excel.ActiveWorkbook

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/840affb00a32037c3048ca1496cc97c5fbf34140/e2e/"

$mdFile1 = "13793503-6295-4ec2-a89c-371c3f162bbe.md"
$mdFile2 = "db62ad9b-b343-4ade-b93f-b5ee8bc325e7.md"

$statusHandedBack = "Handed back: in sync with en-US"

# Column widths observed on the unedited sheets; anything that ends up at
# this COM width keeps its current XML width (no accidental churn).
$wideNarrow = 29.144371396019366   # -> xml col width ~30 (was ~17.22)
$wideFull   = 39.166666666666664   # -> xml col width 40 (was narrower)

# ---------------------------------------------------------------------
# Overview sheet: mark both languages as handed back.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

$overview.Columns.Item(5).ColumnWidth = $wideNarrow
$overview.Columns.Item(6).ColumnWidth = $wideNarrow

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("C3").Value = $statusHandedBack

$zhcn.Range("I2").Value = $mdFile1
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($githubBase + $mdFile1), "", "", $mdFile1) | Out-Null
$zhcn.Range("J2").Value = "13793503-6295-4ec2-a89c-371c3f162bbe.13d87bd803bd70fc72815d62d93ffa80b30ceaff.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-26 15:15:47"

$zhcn.Range("I3").Value = $mdFile2
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($githubBase + $mdFile2), "", "", $mdFile2) | Out-Null
$zhcn.Range("J3").Value = "db62ad9b-b343-4ade-b93f-b5ee8bc325e7.f25a62217b5c7b18cd395b58c1fcba6f71f03f97.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-26 15:15:47"

$zhcn.Columns.Item(3).ColumnWidth = $wideNarrow
$zhcn.Columns.Item(9).ColumnWidth = $wideFull
$zhcn.Columns.Item(10).ColumnWidth = $wideFull

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusHandedBack
$dede.Range("C3").Value = $statusHandedBack

$dede.Range("I2").Value = $mdFile1
$dede.Hyperlinks.Add($dede.Range("I2"), ($githubBase + $mdFile1), "", "", $mdFile1) | Out-Null
$dede.Range("J2").Value = "13793503-6295-4ec2-a89c-371c3f162bbe.13d87bd803bd70fc72815d62d93ffa80b30ceaff.de-de.xlf"
$dede.Range("K2").Value = "2016-08-26 15:15:54"

$dede.Range("I3").Value = $mdFile2
$dede.Hyperlinks.Add($dede.Range("I3"), ($githubBase + $mdFile2), "", "", $mdFile2) | Out-Null
$dede.Range("J3").Value = "db62ad9b-b343-4ade-b93f-b5ee8bc325e7.f25a62217b5c7b18cd395b58c1fcba6f71f03f97.de-de.xlf"
$dede.Range("K3").Value = "2016-08-26 15:15:54"

$dede.Columns.Item(3).ColumnWidth = $wideNarrow
$dede.Columns.Item(9).ColumnWidth = $wideFull
$dede.Columns.Item(10).ColumnWidth = $wideFull

Write-Host "Handback report generated for zh-cn and de-de."
